# Automatic dashboard refresh — "Comentarios" sheet gained a new
# `post_url_original` column, inserted right after `post_url` (between the
# old `post_url` and `author_name` columns). Every existing column from
# `author_name` (E) through `created_time_raw` (N) shifts one column to the
# right (F through O), and the freshly inserted column E is filled with a
# copy of the `post_url` value (column D) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comentarios")

# Insert a blank column before column E ("author_name"); this pushes
# author_name..created_time_raw (E:N) one column to the right (F:O) and
# carries each cell's formatting along with it, matching how Excel's
# "Insert Sheet Columns" command behaves.
$ws.Columns("E:E").Insert()

# Title the newly inserted column.
$ws.Range("E1").Value = "post_url_original"

# Find the last used data row (row 1 is the header).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Backfill the new column with the post_url (column D) for each comment row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 4).Value()
}
